# Config.xlsx update: Start System1_DownloadReports, Navigate to page, enter TaxID,
# select year and select each month.

$wb = $excel.ActiveWorkbook

$wsSettings  = $wb.Worksheets.Item("Settings")
$wsConstants = $wb.Worksheets.Item("Constants")

# --- Settings sheet: update existing values ---
$wsSettings.Range("B2").Value  = "GenerateYearlyReports`n"
$wsSettings.Range("B5").Value  = "FinanceAndAccounting-REF-GenerateYearlyReports-Performer"

# --- Settings sheet: fill in the previously-empty rows 6-13 ---
$wsSettings.Range("A6").Value  = "ReportDirectory"
$wsSettings.Range("B6").Value  = "C:\Users\james.coker\Documents\Reports"

$wsSettings.Range("A7").Value  = "System1_URL"
$wsSettings.Range("B7").Value  = "https://acme-test.uipath.com"

$wsSettings.Range("A8").Value  = "System1_WorkItemsURL"
$wsSettings.Range("B8").Value  = "https://acme-test.uipath.com/work-items/"

$wsSettings.Range("A9").Value  = "System1_DownloardReportURL"
$wsSettings.Range("B9").Value  = "https://acme-test.uipath.com/reports/download"

$wsSettings.Range("A10").Value = "System1_YearlyReportURL"
$wsSettings.Range("B10").Value = "https://acme-test.uipath.com/reports/upload"

$wsSettings.Range("A11").Value = "System1_Credentials"
$wsSettings.Range("B11").Value = "ACMELogin"

$wsSettings.Range("A12").Value = "Status"
$wsSettings.Range("B12").Value = "Completed"

$wsSettings.Range("A13").Value = "Year"
$wsSettings.Range("B13").Value = 2021

# --- Hyperlinks for the URL cells (added in this order to match rId1..rId4) ---
$wsSettings.Hyperlinks.Add($wsSettings.Range("B8"), "https://acme-test.uipath.com/work-items/") | Out-Null
$wsSettings.Hyperlinks.Add($wsSettings.Range("B9"), "https://acme-test.uipath.com/reports/download") | Out-Null
$wsSettings.Hyperlinks.Add($wsSettings.Range("B10"), "https://acme-test.uipath.com/reports/upload") | Out-Null
$wsSettings.Hyperlinks.Add($wsSettings.Range("B7"), "https://acme-test.uipath.com") | Out-Null

# --- Selections / active sheet ---
$wsSettings.Range("B6").Select()
$wsConstants.Activate()
$wsConstants.Range("A12").Select()

Write-Output "edit applied"
